# Update the Load Forecast sheet with the new RES Hourly Production Forecast data.
# Dates shift forward by 8 days (2024-09-10 -> 2024-09-18) and the Load values
# are replaced with the newly forecast figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(45561,45561.01041666666,45561.02083333334,45561.03125,45561.04166666666,45561.05208333334,45561.0625,45561.07291666666,45561.08333333334,45561.09375,45561.10416666666,45561.11458333334,45561.125,45561.13541666666,45561.14583333334,45561.15625,45561.16666666666,45561.17708333334,45561.1875,45561.19791666666,45561.20833333334,45561.21875,45561.22916666666,45561.23958333334,45561.25,45561.26041666666,45561.27083333334,45561.28125,45561.29166666666,45561.30208333334,45561.3125,45561.32291666666,45561.33333333334,45561.34375,45561.35416666666,45561.36458333334,45561.375,45561.38541666666,45561.39583333334,45561.40625,45561.41666666666,45561.42708333334,45561.4375,45561.44791666666,45561.45833333334,45561.46875,45561.47916666666,45561.48958333334,45561.5,45561.51041666666,45561.52083333334,45561.53125,45561.54166666666,45561.55208333334,45561.5625,45561.57291666666,45561.58333333334,45561.59375,45561.60416666666,45561.61458333334,45561.625,45561.63541666666,45561.64583333334,45561.65625,45561.66666666666,45561.67708333334,45561.6875,45561.69791666666,45561.70833333334,45561.71875,45561.72916666666,45561.73958333334,45561.75,45561.76041666666,45561.77083333334,45561.78125,45561.79166666666,45561.80208333334,45561.8125,45561.82291666666,45561.83333333334,45561.84375,45561.85416666666,45561.86458333334,45561.875,45561.88541666666,45561.89583333334,45561.90625,45561.91666666666,45561.92708333334,45561.9375,45561.94791666666)

$loads = @(5080,5040,5000,4960,4910,4880,4860,4860,4850,4850,4870,4880,4900,4920,4940,4970,5010,5080,5170,5290,5410,5550,5690,5840,5970,6090,6200,6280,6350,6390,6400,6400,6390,6360,6300,6230,6150,6060,5980,5900,5830,5770,5710,5670,5640,5610,5580,5570,5550,5540,5540,5550,5560,5590,5620,5660,5690,5730,5770,5810,5850,5900,5950,6010,6080,6160,6240,6320,6400,6480,6560,6660,6760,6870,6970,7060,7150,7200,7200,7170,7050,6900,6770,6620,6440,6300,6160,6010,5860,5720,5600,5480)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]
    $ws.Cells.Item($row, 2).Value = $loads[$i]
}
